$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "feet"
$ws.Range("B4").Value = "ft"
$ws.Range("A5").Value = "meter"
$ws.Range("B5").Value = "m"
$ws.Range("A6").Value = "goal_height_req"
$ws.Range("B6").Value = "Height:"

$ws.Range("B6").Select()
